$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataFeed")

# Update the values to reflect the new data set
$ws.Range("B3").Value = "Mobile"
$ws.Range("B2").Value = "IPAD"

# Move the active selection to B3, matching the new sheet view state
$ws.Range("B3").Select()
